$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update zoom: remove custom zoomScale (reset to 100) and set zoomScaleNormal to 100
$ws.Application.ActiveWindow.Zoom = 100

# New data values (Cilium_number, horizontal_length_in_pixels, z_lower, z_upper)
$data = @(
    @(1, 6.0389999999999997, 2, 10),
    @(2, 6.5, 2, 9),
    @(3, 5.4909999999999997, 2, 10),
    @(4, 9, 3, 9),
    @(5, 7.601, 3, 10),
    @(6, 3.0409999999999999, 6, 12),
    @(7, 7.3109999999999999, 3, 11),
    @(8, 5.1989999999999998, 4, 10),
    @(9, 6.9390000000000001, 3, 11),
    @(10, 9.0879999999999992, 3, 11),
    @(11, 4.8029999999999999, 4, 11),
    @(12, 8.2520000000000007, 4, 11),
    @(13, 9.5879999999999992, 3, 10),
    @(14, 10.272, 4, 12),
    @(15, 10.808, 2, 10),
    @(16, 8.9359999999999999, 4, 11)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

$wb.Save()
